$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) updates in column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 131
$ws1.Range("F9").Value = 2060
$ws1.Range("F10").Value = 360
$ws1.Range("F11").Value = 4933
$ws1.Range("F12").Value = 99
$ws1.Range("F13").Value = 345

# Sheet "全部类型" (All types) updates - same events mirrored, column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 131
$ws4.Range("F13").Value = 2060
$ws4.Range("F14").Value = 360
$ws4.Range("F15").Value = 4933
$ws4.Range("F16").Value = 99
$ws4.Range("F17").Value = 345
